$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B15: change from text "4" to a real number 4
$ws.Range("B15").Value = 4

# Add new row 16 with review data for parisk
$ws.Range("A16").Value = "parisk"
$ws.Range("B16").Value = "'3"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "nan"
$ws.Range("D16").Value = "DIS"
$ws.Range("E16").Value = "WRI"
$ws.Range("F16").Value = "4cbdf296-0ef7-4a60-9d08-bf70fb941ab3"
$ws.Range("G16").Value = "SJTB5GZCb_annotated.xlsx"
$ws.Range("H16").Value = "The paper does not sufficiently discuss and compare the relevant neuroscience literature and related work."
$ws.Range("I16").Value = "Correct"
